# Auto-generated edit script for recommandations.xlsx BRVM update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": rows 2-22, update D (Variation Totale) and E (Derniere Variation) only ---
$ws.Range("D2").Value = 2461.63
$ws.Range("E2").Value = 108.22
$ws.Range("D3").Value = 2045
$ws.Range("E3").Value = 695
$ws.Range("D4").Value = 1970
$ws.Range("E4").Value = 650
$ws.Range("D5").Value = 1944.99
$ws.Range("E5").Value = 655.78
$ws.Range("E6").Value = 590
$ws.Range("E7").Value = 585
$ws.Range("D8").Value = 1745
$ws.Range("E8").Value = 580
$ws.Range("D9").Value = 1570
$ws.Range("E9").Value = 525
$ws.Range("D10").Value = 1068.62
$ws.Range("E10").Value = 357.9
$ws.Range("E11").Value = 345.13
$ws.Range("D12").Value = 923.83
$ws.Range("E12").Value = 308.61
$ws.Range("D13").Value = 707.87
$ws.Range("E13").Value = 240.84
$ws.Range("D14").Value = 588.27
$ws.Range("E14").Value = 200.08
$ws.Range("D15").Value = 549.95
$ws.Range("E15").Value = 183.87
$ws.Range("D16").Value = 412.6
$ws.Range("E16").Value = 134.55
$ws.Range("D17").Value = 389.87
$ws.Range("E17").Value = 130.05
$ws.Range("D18").Value = 369.32
$ws.Range("E18").Value = 122.81
$ws.Range("D19").Value = 362.96
$ws.Range("E19").Value = 120.7
$ws.Range("D20").Value = 317.49
$ws.Range("E20").Value = 106.44
$ws.Range("D21").Value = 312.79
$ws.Range("E21").Value = 104.65
$ws.Range("D22").Value = 277.68
$ws.Range("E22").Value = 92.17

# --- Sheet "Recommandations": rows 24-44, full content re-sorted by Variation Totale (%) ---
$ws.Range("A24").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 12.87
$ws.Range("E24").Value = 7.4
$ws.Range("F24").Value = "🟡 Observer"
$ws.Range("G24").Value = "➖ Neutre"

$ws.Range("A25").Value = "BERNABE CI (BNBC)"
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 9.67
$ws.Range("E25").Value = 4.71
$ws.Range("F25").Value = "🟡 Observer"
$ws.Range("G25").Value = "👀 À surveiller"

$ws.Range("A26").Value = "BICI CI (BICC)"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 7.48
$ws.Range("E26").Value = 7.48
$ws.Range("F26").Value = "🟡 Observer"
$ws.Range("G26").Value = "➖ Neutre"

$ws.Range("A27").Value = "FILTISAC CI (FTSC)"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 5.93
$ws.Range("E27").Value = -3.5
$ws.Range("F27").Value = "🟡 Observer"
$ws.Range("G27").Value = "👀 À surveiller"

$ws.Range("A28").Value = "ORAGROUP TOGO (ORGT)"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 5.83
$ws.Range("E28").Value = 5.83
$ws.Range("F28").Value = "🟡 Observer"
$ws.Range("G28").Value = "➖ Neutre"

$ws.Range("A29").Value = "SOLIBRA CI (SLBC)"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 3.83
$ws.Range("E29").Value = 3.83
$ws.Range("F29").Value = "🟡 Observer"
$ws.Range("G29").Value = "➖ Neutre"

$ws.Range("A30").Value = "SAFCA CI (SAFC)"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 2.99
$ws.Range("E30").Value = 2.99
$ws.Range("F30").Value = "🟡 Observer"
$ws.Range("G30").Value = "➖ Neutre"

$ws.Range("A31").Value = "UNIWAX CI (UNXC)"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2.1
$ws.Range("E31").Value = 7.27
$ws.Range("F31").Value = "🟡 Observer"
$ws.Range("G31").Value = "👀 À surveiller"

$ws.Range("A32").Value = "SOGB CI (SOGC)"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 1.8
$ws.Range("E32").Value = 1.8
$ws.Range("F32").Value = "🟡 Observer"
$ws.Range("G32").Value = "➖ Neutre"

$ws.Range("A33").Value = "TOTAL"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "🟡 Observer"
$ws.Range("G33").Value = "➖ Neutre"

$ws.Range("A34").Value = "NEI-CEDA CI (NEIC)"
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = -1.68
$ws.Range("E34").Value = -1.68
$ws.Range("F34").Value = "🟡 Observer"
$ws.Range("G34").Value = "➖ Neutre"

$ws.Range("A35").Value = "SETAO CI (STAC)"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = -1.71
$ws.Range("E35").Value = -1.71
$ws.Range("F35").Value = "🟡 Observer"
$ws.Range("G35").Value = "➖ Neutre"

$ws.Range("A36").Value = "AIR LIQUIDE CI (SIVC)"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = -1.92
$ws.Range("E36").Value = -1.92
$ws.Range("F36").Value = "🟡 Observer"
$ws.Range("G36").Value = "➖ Neutre"

$ws.Range("A37").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = -2.02
$ws.Range("E37").Value = -2.02
$ws.Range("F37").Value = "🟡 Observer"
$ws.Range("G37").Value = "➖ Neutre"

$ws.Range("A38").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = -2.42
$ws.Range("E38").Value = -2.42
$ws.Range("F38").Value = "🟡 Observer"
$ws.Range("G38").Value = "➖ Neutre"

$ws.Range("A39").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = -2.51
$ws.Range("E39").Value = -2.51
$ws.Range("F39").Value = "🟡 Observer"
$ws.Range("G39").Value = "➖ Neutre"

$ws.Range("A40").Value = "SONATEL SN (SNTS)"
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = -3.81
$ws.Range("E40").Value = -3.81
$ws.Range("F40").Value = "🟡 Observer"
$ws.Range("G40").Value = "➖ Neutre"

$ws.Range("A41").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = -4.26
$ws.Range("E41").Value = -4.26
$ws.Range("F41").Value = "🟡 Observer"
$ws.Range("G41").Value = "➖ Neutre"

$ws.Range("A42").Value = "ONATEL BF (ONTBF)"
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = -4.37
$ws.Range("E42").Value = -4.37
$ws.Range("F42").Value = "🟡 Observer"
$ws.Range("G42").Value = "➖ Neutre"

$ws.Range("A43").Value = "BANK OF AFRICA NG (BOAN)"
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = -5.05
$ws.Range("E43").Value = -5.05
$ws.Range("F43").Value = "🟡 Observer"
$ws.Range("G43").Value = "➖ Neutre"

$ws.Range("A44").Value = "BANK OF AFRICA BN (BOAB)"
$ws.Range("B44").Value = 0
$ws.Range("C44").Value = 2
$ws.Range("D44").Value = -5.89
$ws.Range("E44").Value = -2.44
$ws.Range("F44").Value = "🟡 Observer"
$ws.Range("G44").Value = "➖ Neutre"

# Row 45 (ONATEL BF) no longer exists as a separate row - remove it, shifting dimension to A1:G44
$ws.Rows.Item(45).Delete()

# --- Sheet "Top_YTD": update Progression YTD (%) values ---
$ws2.Range("B2").Value = 489150.58
$ws2.Range("B3").Value = 47649.69
$ws2.Range("B4").Value = 43203.13
$ws2.Range("B5").Value = 41790.19
$ws2.Range("B8").Value = 31574.4
$ws2.Range("B9").Value = 24118.75
$ws2.Range("B10").Value = 9394.68
